$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '44.090.89'
$ws.Range("E2").Value = '  +0.85%  '

$ws.Range("D3").Value = '2.366.96'
$ws.Range("E3").Value = '  +0.69%  '

$ws.Range("E4").Value = '  +0.31%  '

$ws.Range("D5").Value = '''0.697'
$ws.Range("E5").Value = '  +5.83%  '

$ws.Range("D6").Value = '''241.77'
$ws.Range("E6").Value = '  +3.03%  '

$ws.Range("D7").Value = '''76.96'
$ws.Range("E7").Value = '  +5.22%  '

$ws.Range("E8").Value = '  +0.00%  '

$ws.Range("D9").Value = '''0.623'
$ws.Range("E9").Value = '  +16.47%  '

$ws.Range("E10").Value = '  +3.55%  '

$ws.Range("D11").Value = '''57.54'
$ws.Range("E11").Value = '  +1.19%  '

$ws.Range("D12").Value = '''33.62'
$ws.Range("E12").Value = '  +18.66%  '

$ws.Range("E13").Value = '  +13.48%  '

$ws.Range("E14").Value = '  +2.31%  '

$ws.Range("D15").Value = '2.720.42'
$ws.Range("E15").Value = '  +0.61%  '

$ws.Range("D16").Value = '''16.75'
$ws.Range("E16").Value = '  -0.62%  '

$ws.Range("D17").Value = '''0.931'
$ws.Range("E17").Value = '  +5.26%  '

$ws.Range("D18").Value = '2.353.05'
$ws.Range("E18").Value = '  -0.11%  '

$ws.Range("D19").Value = '44.038.46'
$ws.Range("E19").Value = '  +1.03%  '

$ws.Range("E20").Value = '  +2.36%  '

$ws.Range("D21").Value = '''6.73'
$ws.Range("E21").Value = '  +6.64%  '

$ws.Range("D22").Value = '''77.83'
$ws.Range("E22").Value = '  +2.52%  '

$ws.Range("D23").Value = '''259.77'
$ws.Range("E23").Value = '  +3.67%  '

$ws.Range("E24").Value = '  +0.24%  '

$ws.Range("D25").Value = '''3.74'
$ws.Range("E25").Value = '  -1.73%  '

$ws.Range("E26").Value = '  +2.69%  '

$ws.Range("E27").Value = '  +16.86%  '

$ws.Range("D28").Value = '''10.92'
$ws.Range("E28").Value = '  +6.43%  '

$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").Value = '''2.29'
$ws.Range("E29").Value = '  +1.94%  '

$ws.Range("B30").Value = 'EthereumClassic'
$ws.Range("C30").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D30").Value = '''23.26'
$ws.Range("E30").Value = '  +3.89%  '

$ws.Range("D31").Value = '''175.32'
$ws.Range("E31").Value = '  +1.70%  '

$ws.Range("E32").Value = '  -3.09%  '

$ws.Range("E33").Value = '  +4.95%  '

$ws.Range("E34").Value = '  +4.22%  '

$ws.Range("D35").Value = '''0.0768'
$ws.Range("E35").Value = '  +9.57%  '

$ws.Range("D36").Value = '''5.44'
$ws.Range("E36").Value = '  +5.97%  '

$ws.Range("D37").Value = '''3.78'
$ws.Range("E37").Value = '  +0.81%  '

$ws.Range("D38").Value = '''2.43'
$ws.Range("E38").Value = '  +0.32%  '

$ws.Range("E39").Value = '  +1.28%  '

$ws.Range("E40").Value = '  +7.45%  '

$ws.Range("D41").Value = '''0.221'
$ws.Range("E41").Value = '  +23.07%  '

$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").Value = '''9.22'
$ws.Range("E42").Value = '  +3.84%  '

$ws.Range("B43").Value = 'InjectiveProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D43").Value = '''19.35'
$ws.Range("E43").Value = '  +0.21%  '

$ws.Range("D44").Value = '''0.107'
$ws.Range("E44").Value = '  +11.89%  '

$ws.Range("D45").Value = '''4.89'
$ws.Range("E45").Value = '  +10.29%  '

$ws.Range("E46").Value = '  -0.08%  '

$ws.Range("E47").Value = '  +11.41%  '

$ws.Range("B48").Value = 'Aave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D48").Value = '''102.99'
$ws.Range("E48").Value = '  +4.09%  '

$ws.Range("B49").Value = 'TrustWalletToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D49").Value = '''1.26'
$ws.Range("E49").Value = '  +4.28%  '

$ws.Range("B50").Value = 'ARBITRUM'
$ws.Range("C50").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D50").Value = '''1.20'
$ws.Range("E50").Value = '  +2.88%  '

$ws.Range("D51").Value = '''56.14'
$ws.Range("E51").Value = '  +8.59%  '
